# THNV3380EP105 - September - Letter 2
# 1) Bump the letter date forward two days.
# 2) Split the manager's mailing address onto its own line (street / city-state-zip)
#    and leave an extra blank line afterwards.
# 3) Drop the two blank paragraphs that used to sit between "Board of Directors"
#    and the Fair Debt Collection Practices Act notice.

$d = $word.ActiveDocument

# --- 1) Remove the two blank paragraphs right after "Board of Directors" ----------
# Do this first, while paragraph indices still match the pristine document, so the
# later Find/Replace calls (which are text-anchored, not index-anchored) aren't
# affected by the reindex.
$boardParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Board of Directors*") {
        $boardParaIndex = $i
        break
    }
}

if ($boardParaIndex -gt 0) {
    $pBlank1 = $d.Paragraphs.Item($boardParaIndex + 1)
    $pBlank2 = $d.Paragraphs.Item($boardParaIndex + 2)
    $delRange = $d.Range($pBlank1.Range.Start, $pBlank2.Range.End)
    # Each Delete() call only consumes one paragraph mark in this host, so call it
    # once per paragraph being removed.
    $delRange.Delete()
    $delRange.Delete()
}

# --- 2) Update the letter date -----------------------------------------------------
$d.Content.Find.Execute("September 19, 2025", $false, $false, $false, $false, $false, `
    $true, 1, $false, "September 21, 2025", 2)

# --- 3) Split "1550 Technology Dr., San Jose CA 95110" into two lines, plus a ------
#        trailing blank line, matching the new mailing-address block.
$d.Content.Find.Execute("1550 Technology Dr., San Jose CA 95110", $false, $false, $false, $false, $false, `
    $true, 1, $false, "1550 Technology Dr.^pSan Jose, CA 95110^p", 2)
